$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Competitions")

# New competition rows for the last 6 months of 2022.
# Entry order matters for shared-string allocation order (matches the
# original author's edit order: rows 51-55 first, then row 50).
$newRows = @(
    @{ Row = 51; Date = 44782; Style = "Club Yeast Beers";                     Year = 2015 },
    @{ Row = 52; Date = 44817; Style = "Meads, Cysers, Ciders";                Year = 2015 },
    @{ Row = 53; Date = 44845; Style = "Open";                                 Year = 2015 },
    @{ Row = 54; Date = 44873; Style = "SHIVs, Darks (NOT Porters or Stouts)"; Year = 2015 },
    @{ Row = 55; Date = 44908; Style = "Porters, Stouts, Big Beers";           Year = 2015 },
    @{ Row = 50; Date = 44754; Style = "Session Beers (under 5% ABV)";         Year = 2015 }
)

foreach ($r in $newRows) {
    $ws.Cells.Item($r.Row, 1).Value = $r.Date
    $ws.Cells.Item($r.Row, 2).Value = $r.Style
    $ws.Cells.Item($r.Row, 4).Value = $r.Year
}

# Expand the Competitions table (Table4) and its autofilter to cover the new rows
$lo = $ws.ListObjects.Item(1)
$lo.Resize($ws.Range("A1:D55"))

# Update view/selection state: Competitions becomes the active sheet/tab
$ws.Activate() | Out-Null
$ws.Range("C51").Select() | Out-Null

$wsEntries = $wb.Worksheets.Item("Entries")
$wsEntries.Activate() | Out-Null
$wsEntries.Range("A145").Select() | Out-Null
$winEntries = $excel.ActiveWindow
$winEntries.ScrollRow = 127
$winEntries.ScrollColumn = 1

$ws.Activate() | Out-Null
$win = $excel.ActiveWindow
$win.ScrollRow = 29
$win.ScrollColumn = 1
